$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 26 and 29 had their match data (columns F:V) swapped.
#    Columns A:E (index, pais, torneio, temporada, data_partida) stay as-is.
# ---------------------------------------------------------------------------

$row26 = @("Kluczbork", 3, "Warta Gorzow", 0, 1.85, "22/08/2023 05:12", 2.36, "23/08/2023 16:57", 3.38, "22/08/2023 05:12", 3.09, "23/08/2023 16:56", 3.18, "22/08/2023 05:12", 2.8, "23/08/2023 16:57", "https://www.betexplorer.com/football/poland/iii-liga-group-iii/kluczbork-warta-gorzow/W2oakoGJ/")

$row29 = @("Polkowice", 2, "Starowice Dolne", 2, 1.53, "23/08/2023 11:12", 1.71, "23/08/2023 16:37", 4.29, "23/08/2023 11:12", 4.53, "23/08/2023 16:37", 4.42, "23/08/2023 11:12", 3.31, "23/08/2023 16:37", "https://www.betexplorer.com/football/poland/iii-liga-group-iii/polkowice-starowice-dolne/txyZrqgm/")

for ($i = 0; $i -lt $row26.Length; $i++) {
    $ws.Cells.Item(26, 6 + $i).Value = $row26[$i]
}
for ($i = 0; $i -lt $row29.Length; $i++) {
    $ws.Cells.Item(29, 6 + $i).Value = $row29[$i]
}

# ---------------------------------------------------------------------------
# 2) Append 5 new match rows (63..67), reusing the formatting of row 62 by
#    copying it down first, then overwriting the per-row values.
# ---------------------------------------------------------------------------

$newRows = @(
    @{ A=62; E=45192.58333333334; F="Goczalkowice Zdroj"; G=2; H="Zielona Gora"; I=0;
       J=2.04; K="22/09/2023 01:12"; L=2.34; M="23/09/2023 13:57";
       N=3.18; O="22/09/2023 01:12"; P=3.04; Q="23/09/2023 13:51";
       R=2.91; S="22/09/2023 01:12"; T=2.87; U="23/09/2023 13:57";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iii/goczalkowice-zdroj-zielona-gora/bgi4BPY7/" },
    @{ A=63; E=45192.58333333334; F="Warta Gorzow"; G=1; H="Pawlowice"; I=3;
       J=2.47; K="22/09/2023 01:12"; L=2.99; M="23/09/2023 13:36";
       N=3.1;  O="22/09/2023 01:12"; P=3.45; Q="23/09/2023 13:36";
       R=2.36; S="22/09/2023 01:12"; T=2.08; U="23/09/2023 13:36";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iii/warta-gorzow-pniowek-pawlowice/QaqjZSQE/" },
    @{ A=64; E=45192.66666666666; F="Gwarek Tarnowskie Gory"; G=0; H="Polkowice"; I=1;
       J=2.12; K="22/09/2023 03:13"; L=2.31; M="23/09/2023 15:47";
       N=3.22; O="22/09/2023 03:13"; P=3.42; Q="23/09/2023 15:42";
       R=2.75; S="22/09/2023 03:13"; T=2.63; U="23/09/2023 15:47";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iii/gwarek-tarnowskie-gory-polkowice/nFXbX6eR/" },
    @{ A=65; E=45193.58333333334; F="Gornik Zabrze II"; G=2; H="Slask Wroclaw II"; I=3;
       J=2.83; K="23/09/2023 01:12"; L=3.16; M="24/09/2023 13:51";
       N=3.35; O="23/09/2023 01:12"; P=3.69; Q="24/09/2023 13:51";
       R=2.02; S="23/09/2023 01:12"; T=1.93; U="24/09/2023 13:51";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iii/gornik-zabrze-slask-wroclaw/8OYfYntL/" },
    @{ A=66; E=45193.66666666666; F="Unia Turza Slaska"; G=1; H="Jelenia Gora"; I=1;
       J=2.54; K="24/09/2023 13:42"; L=2.46; M="24/09/2023 15:18";
       N=3.4;  O="24/09/2023 13:42"; P=3.85; Q="24/09/2023 15:18";
       R=2.34; S="24/09/2023 13:42"; T=2.28; U="24/09/2023 15:18";
       V="https://www.betexplorer.com/football/poland/iii-liga-group-iii/unia-turza-slaska-karkonosze-jelenia-gora/2RceDo4e/" }
)

$srcRow = $ws.Range("A62:V62")
$startRow = 63

for ($idx = 0; $idx -lt $newRows.Length; $idx++) {
    $r = $startRow + $idx
    $data = $newRows[$idx]

    $destRow = $ws.Range("A" + $r + ":V" + $r)
    $srcRow.Copy($destRow)

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = "poland"
    $ws.Cells.Item($r, 3).Value = "iii-liga-group-iii"
    $ws.Cells.Item($r, 4).Value = "2023-2024"
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = $data.T
    $ws.Cells.Item($r, 21).Value = $data.U
    $ws.Cells.Item($r, 22).Value = $data.V
}
